$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92 (shifts existing rows 92..146 down to 93..147)
$ws.Rows(92).Insert()

# Populate the newly inserted row 92 with the new weekly record
$ws.Cells.Item(92, 1).Value = 4
$ws.Cells.Item(92, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(92, 3).Value = "Los Lagos"
$ws.Cells.Item(92, 4).Value = 45126
$ws.Cells.Item(92, 5).Value = 10
$ws.Cells.Item(92, 6).Value = 100112026
$ws.Cells.Item(92, 7).Value = "Haba"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 80
$ws.Cells.Item(92, 11).Value = 21000
$ws.Cells.Item(92, 12).Value = 21000
$ws.Cells.Item(92, 13).Value = 21000
$ws.Cells.Item(92, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(92, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(92, 16).Value = 840
$ws.Cells.Item(92, 17).Value = 25
$ws.Cells.Item(92, 18).Value = "Hortaliza"
